$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revise figures for the already-existing months (01-2021 .. 05-2021) ---
# Row 218: 01-01-2021
$ws.Range("B218").Value = 210862
$ws.Range("C218").Value = 36657
$ws.Range("F218").Value = 2097
$ws.Range("G218").Value = 89970
$ws.Range("H218").Value = 56697
$ws.Range("I218").Value = 203821

# Row 219: 01-02-2021
$ws.Range("B219").Value = 209836
$ws.Range("C219").Value = 35564
$ws.Range("F219").Value = 2122
$ws.Range("G219").Value = 89284
$ws.Range("H219").Value = 57221
$ws.Range("I219").Value = 205126

# Row 220: 01-03-2021
$ws.Range("B220").Value = 209510
$ws.Range("C220").Value = 33754
$ws.Range("E220").Value = 24767
$ws.Range("F220").Value = 2125
$ws.Range("G220").Value = 90203
$ws.Range("H220").Value = 57157
$ws.Range("I220").Value = 205567

# Row 221: 01-04-2021
$ws.Range("B221").Value = 213397
$ws.Range("C221").Value = 35862
$ws.Range("E221").Value = 24594
$ws.Range("F221").Value = 2131
$ws.Range("G221").Value = 91320
$ws.Range("H221").Value = 58214
$ws.Range("I221").Value = 209412

# Row 222: 01-05-2021
$ws.Range("B222").Value = 214870
$ws.Range("C222").Value = 37594
$ws.Range("F222").Value = 2100
$ws.Range("G222").Value = 90877
$ws.Range("H222").Value = 57870
$ws.Range("I222").Value = 211503

# --- Append the new month: 01-06-2021 (row 223) ---
# Prefix with an apostrophe so the "dd-mm-yyyy"-like text is kept as a plain
# text label (matching the rest of column A) instead of being auto-converted
# into a date serial number, then reset the style so no extra number format
# sticks to the cell.
$ws.Range("A223").Value = "'01-06-2021"
$ws.Range("A223").Style = "Normal"

$ws.Range("B223").Value = 213283
$ws.Range("C223").Value = 37860
$ws.Range("D223").Value = 1265
$ws.Range("E223").Value = 24276
$ws.Range("F223").Value = 2081
$ws.Range("G223").Value = 91116
$ws.Range("H223").Value = 56685
$ws.Range("I223").Value = 209472
